$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 436
$ws.Range("I111").Value = 450
$ws.Range("J111").Value = 432.5
$ws.Range("K111").Value = 1350
$ws.Range("L111").Value = 1297.5
$ws.Range("M111").Value = 1717
$ws.Range("N111").Value = -7431.5
$ws.Range("H125").Value = 1628.3636
$ws.Range("I125").Value = 1266.4
$ws.Range("J125").Value = 1930
$ws.Range("K125").Value = 11397.6
$ws.Range("L125").Value = 17370
$ws.Range("M125").Value = -8937.6
$ws.Range("N125").Value = -22290

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 128
$ws.Range("I5").Value = 144.5
$ws.Range("J5").Value = 95
$ws.Range("K5").Value = 144.5
$ws.Range("L5").Value = 95
$ws.Range("M5").Value = -32.5
$ws.Range("N5").Value = -319
$ws.Range("H21").Value = 60000
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").ClearContents()
$ws.Range("H30").Value = 35250
$ws.Range("I30").Value = 500
$ws.Range("J30").Value = 70000
$ws.Range("K30").Value = 500
$ws.Range("L30").Value = 70000
$ws.Range("M30").Value = -350
$ws.Range("N30").Value = -70300
$ws.Range("H32").Value = 22475.78
$ws.Range("I32").Value = 3819.103
$ws.Range("K32").Value = 3819.103
$ws.Range("M32").Value = -3532.103
$ws.Range("H33").Value = 62500
$ws.Range("I33").Value = 5000
$ws.Range("K33").Value = 5000
$ws.Range("M33").Value = -4671
$ws.Range("H36").Value = 43463
$ws.Range("I36").Value = 4926
$ws.Range("J36").Value = 82000
$ws.Range("K36").Value = 4926
$ws.Range("L36").Value = 82000
$ws.Range("M36").Value = -4580
$ws.Range("N36").Value = -82692
$ws.Range("H45").Value = 1401.5172
$ws.Range("I45").Value = 1357.1852
$ws.Range("K45").Value = 1357.1852
$ws.Range("M45").Value = -980.1851999999999
$ws.Range("H110").Value = 896.92
$ws.Range("I110").Value = 710.7222
$ws.Range("J110").Value = 1375.7142
$ws.Range("K110").Value = 710.7222
$ws.Range("L110").Value = 1375.7142
$ws.Range("M110").Value = 1334.2778
$ws.Range("N110").Value = -5465.7142

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 128
$ws.Range("I4").Value = 144.5
$ws.Range("J4").Value = 95
$ws.Range("K4").Value = 144.5
$ws.Range("L4").Value = 95
$ws.Range("M4").Value = -29.5
$ws.Range("N4").Value = -325
$ws.Range("H22").Value = 1429.2858
$ws.Range("I22").Value = 52.5
$ws.Range("J22").Value = 1980
$ws.Range("K22").Value = 52.5
$ws.Range("L22").Value = 1980
$ws.Range("M22").Value = 120.5
$ws.Range("N22").Value = -2326
$ws.Range("H134").Value = 1513.0286
$ws.Range("I134").Value = 1428.3793
$ws.Range("J134").Value = 1922.1666
$ws.Range("K134").Value = 4285.1379
$ws.Range("L134").Value = 5766.4998
$ws.Range("M134").Value = -1750.1379
$ws.Range("N134").Value = -10836.4998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H32").Value = 5000
$ws.Range("I32").Value = 5000
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 5000
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -4684
$ws.Range("N32").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1958
$ws.Range("J113").Value = 497.5
$ws.Range("L113").Value = 1492.5
$ws.Range("N113").Value = -5832.5
$ws.Range("H131").Value = 637.5
$ws.Range("I131").Value = 585.7143
$ws.Range("J131").Value = 1000
$ws.Range("K131").Value = 1757.1429
$ws.Range("L131").Value = 3000
$ws.Range("M131").Value = 3282.8571
$ws.Range("N131").Value = -13080

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 55000
$ws.Range("J18").Value = 55000
$ws.Range("L18").Value = 55000
$ws.Range("N18").Value = -55586
$ws.Range("H21").Value = 16705000
$ws.Range("I21").Value = 50000000
$ws.Range("J21").Value = 57500
$ws.Range("K21").Value = 50000000
$ws.Range("L21").Value = 57500
$ws.Range("M21").Value = -49999827
$ws.Range("N21").Value = -57846
$ws.Range("H30").Value = 16705000
$ws.Range("I30").Value = 50000000
$ws.Range("J30").Value = 57500
$ws.Range("K30").Value = 50000000
$ws.Range("L30").Value = 57500
$ws.Range("M30").Value = -49999895
$ws.Range("N30").Value = -57710

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("M3").ClearContents()
$ws.Range("N3").ClearContents()
$ws.Range("H7").Value = 2217.0908
$ws.Range("I7").Value = 1763
$ws.Range("J7").Value = 2762
$ws.Range("K7").Value = 1763
$ws.Range("L7").Value = 2762
$ws.Range("M7").Value = -1651
$ws.Range("N7").Value = -2986
$ws.Range("H13").Value = 4165
$ws.Range("J13").Value = 4165
$ws.Range("L13").Value = 4165
$ws.Range("N13").Value = -4445
$ws.Range("H14").Value = 7621.273
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 7621.273
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 7621.273
$ws.Range("M14").ClearContents()
$ws.Range("N14").Value = -7965.273
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("M15").ClearContents()
$ws.Range("N15").ClearContents()
$ws.Range("H22").Value = 550.25
$ws.Range("I22").Value = 1001
$ws.Range("J22").Value = 400
$ws.Range("K22").Value = 1001
$ws.Range("L22").Value = 400
$ws.Range("M22").Value = -706
$ws.Range("N22").Value = -990
$ws.Range("H26").Value = 35000
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 35000
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 35000
$ws.Range("M26").ClearContents()
$ws.Range("N26").Value = -35590
$ws.Range("H27").Value = 550.25
$ws.Range("I27").Value = 1001
$ws.Range("J27").Value = 400
$ws.Range("K27").Value = 1001
$ws.Range("L27").Value = 400
$ws.Range("M27").Value = -894
$ws.Range("N27").Value = -614
$ws.Range("H29").Value = 43333.332
$ws.Range("I29").Value = 30000
$ws.Range("J29").Value = 70000
$ws.Range("K29").Value = 30000
$ws.Range("L29").Value = 70000
$ws.Range("M29").Value = -29705
$ws.Range("N29").Value = -70590
$ws.Range("H31").Value = 7929.125
$ws.Range("I31").Value = 833.3333
$ws.Range("J31").Value = 12186.6
$ws.Range("K31").Value = 833.3333
$ws.Range("L31").Value = 12186.6
$ws.Range("M31").Value = -585.3333
$ws.Range("N31").Value = -12682.6
$ws.Range("H34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("N34").ClearContents()
$ws.Range("H126").Value = 2217.0908
$ws.Range("I126").Value = 1763
$ws.Range("J126").Value = 2762
$ws.Range("K126").Value = 5289
$ws.Range("L126").Value = 8286
$ws.Range("M126").Value = -2819
$ws.Range("N126").Value = -13226

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 5000
$ws.Range("J14").Value = 5000
$ws.Range("L14").Value = 5000
$ws.Range("N14").Value = -5336
$ws.Range("H29").Value = 19500
$ws.Range("I29").Value = 8000
$ws.Range("J29").Value = 23333.334
$ws.Range("K29").Value = 8000
$ws.Range("L29").Value = 23333.334
$ws.Range("M29").Value = -7710
$ws.Range("N29").Value = -23913.334
$ws.Range("H32").Value = 20131.5
$ws.Range("I32").Value = 3508.6667
$ws.Range("J32").Value = 70000
$ws.Range("K32").Value = 3508.6667
$ws.Range("L32").Value = 70000
$ws.Range("M32").Value = -3191.6667
$ws.Range("N32").Value = -70634
$ws.Range("H34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("N34").ClearContents()
$ws.Range("H136").Value = 737.1739
$ws.Range("I136").Value = 578.8889
$ws.Range("J136").Value = 1307
$ws.Range("K136").Value = 1736.6667
$ws.Range("L136").Value = 3921
$ws.Range("M136").Value = 813.3332999999998
$ws.Range("N136").Value = -9021
